# Update the main GSC export data on the "Chart" sheet:
#  - Drop the oldest date row (2025-10-13); every later row shifts up one.
#  - Append a new row for the newest date (2026-01-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the first data row (row 2, which holds 2025-10-13). This shifts
# every following row up by one, matching the diff where each row's
# "HTTPS URLs" value becomes the value that used to belong to the next row.
$ws.Rows.Item(2).Delete()

# After the delete, the used range now ends at row 90 (2026-01-10).
# Append the new last row for 2026-01-11. Pre-format column A on this row
# as Text so the date-formatted string isn't auto-converted into a date
# serial number (matching the source data, which stores dates as plain text).
$newRow = 91
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2026-01-11"
$ws.Cells.Item($newRow, 2).Value = 0.0
$ws.Cells.Item($newRow, 3).Value = 26.0

$wb.Save()
